$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from pt-PT to pt-BR
$ws.Name = "pt-BR"

# Update the header text in F2 (table column header "Portuguese (pt-PT)" -> "Portuguese (pt-BR)")
$ws.Range("F2").Value = "Portuguese (pt-BR)"
